$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: mark task complete (checkbox-style boolean in F6)
$ws.Range("F6").Value = $true

# Row 7: fill in actual/predicted start + actual end dates, mark complete
$ws.Range("B7").Value = "5/1/2025"
$ws.Range("C7").Value = "5/1/2025"
$ws.Range("D7").Value = "5/6/2025"
$ws.Range("F7").Value = $true

# Row 8: note the predicted start date as free text
$ws.Range("B8").Value = "5/62025"

# Leave the final selection on F6, matching the saved sheet view
$ws.Range("F6").Select() | Out-Null
